# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: update title (D5) and link (E5)
$ws.Range("D5").Value = "체르노프 유계(Chernoff Bound)"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/09/13/Chernoff_Bound.html"

# Row 25: update title (D25) and link (E25)
$ws.Range("D25").Value = "[바람돌이/딥러닝] Speech - 음성 데이터 이론 및 이해"
$ws.Range("E25").Value = "https://blog.naver.com/winddori2002/222872853715"

# Row 27: update title (D27) and link (E27)
$ws.Range("D27").Value = "TFX 머신러닝 파이프라인 사용하기"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/use-tfx-pipeline-with-customization/"
